$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
# ALC row 33
$ws1.Range("H33").Value = 45335.234
$ws1.Range("J33").Value = 4064
$ws1.Range("L33").Value = 4064
$ws1.Range("N33").Value = -4522
# ALC row 64
$ws1.Range("H64").Value = 3590.818
$ws1.Range("I64").Value = 3000
$ws1.Range("J64").Value = 3812.375
$ws1.Range("K64").Value = 3000
$ws1.Range("L64").Value = 3812.375
$ws1.Range("M64").Value = -2752
$ws1.Range("N64").Value = -4308.375
# ALC row 67
$ws1.Range("H67").Value = 3590.818
$ws1.Range("I67").Value = 3000
$ws1.Range("J67").Value = 3812.375
$ws1.Range("K67").Value = 3000
$ws1.Range("L67").Value = 3812.375
$ws1.Range("M67").Value = -2142
$ws1.Range("N67").Value = -5528.375
# ALC row 69
$ws1.Range("H69").Value = 6185.7144
$ws1.Range("I69").Value = 0
$ws1.Range("J69").Value = 6185.7144
$ws1.Range("K69").Value = 0
$ws1.Range("L69").Value = 18557.1432
$ws1.Range("M69").ClearContents()
$ws1.Range("N69").Value = -20305.1432
# ALC row 72
$ws1.Range("H72").Value = 6185.7144
$ws1.Range("I72").Value = 0
$ws1.Range("J72").Value = 6185.7144
$ws1.Range("K72").Value = 0
$ws1.Range("L72").Value = 55671.4296
$ws1.Range("M72").ClearContents()
$ws1.Range("N72").Value = -64407.4296
# ALC row 76
$ws1.Range("H76").Value = 4965749.5
$ws1.Range("I76").Value = 4237.875
$ws1.Range("K76").Value = 4237.875
$ws1.Range("M76").Value = -3922.875
# ALC row 79
$ws1.Range("H79").Value = 4965749.5
$ws1.Range("I79").Value = 4237.875
$ws1.Range("K79").Value = 4237.875
$ws1.Range("M79").Value = -3145.875
# ALC row 127
$ws1.Range("H127").Value = 1764.7059
$ws1.Range("I127").Value = 1561.6923
$ws1.Range("K127").Value = 4685.0769
$ws1.Range("M127").Value = 274.9231
# ALC row 137
$ws1.Range("H137").Value = 2617.7856
$ws1.Range("I137").Value = 2129.7058
$ws1.Range("K137").Value = 6389.117400000001
$ws1.Range("M137").Value = -3839.117400000001

$ws2 = $wb.Worksheets.Item("ARM")
# ARM row 17
$ws2.Range("H17").Value = 2755
$ws2.Range("I17").Value = 2504
$ws2.Range("J17").Value = 3006
$ws2.Range("K17").Value = 2504
$ws2.Range("L17").Value = 3006
$ws2.Range("M17").Value = -2331
$ws2.Range("N17").Value = -3352
# ARM row 37
$ws2.Range("H37").Value = 12000
$ws2.Range("J37").Value = 12000
$ws2.Range("L37").Value = 12000
$ws2.Range("N37").Value = -12546
# ARM row 45
$ws2.Range("H45").Value = 2335.5
$ws2.Range("I45").Value = 1499.75
$ws2.Range("K45").Value = 1499.75
$ws2.Range("M45").Value = -1122.75
# ARM row 88
$ws2.Range("H88").Value = 12822048
$ws2.Range("I88").Value = 33334284
$ws2.Range("J88").Value = 1900.125
$ws2.Range("K88").Value = 33334284
$ws2.Range("L88").Value = 1900.125
$ws2.Range("M88").Value = -33333878
$ws2.Range("N88").Value = -2712.125
# ARM row 91
$ws2.Range("H91").Value = 12822048
$ws2.Range("I91").Value = 33334284
$ws2.Range("J91").Value = 1900.125
$ws2.Range("K91").Value = 33334284
$ws2.Range("L91").Value = 1900.125
$ws2.Range("M91").Value = -33332880
$ws2.Range("N91").Value = -4708.125
# ARM row 97
$ws2.Range("H97").Value = 1159.7241
$ws2.Range("I97").Value = 1290.9
$ws2.Range("J97").Value = 868.2222
$ws2.Range("K97").Value = 1290.9
$ws2.Range("L97").Value = 868.2222
$ws2.Range("M97").Value = -794.9000000000001
$ws2.Range("N97").Value = -1860.2222

$ws3 = $wb.Worksheets.Item("BSM")
# BSM row 11
$ws3.Range("H11").Value = 1054.2858
$ws3.Range("J11").Value = 1396
$ws3.Range("L11").Value = 1396
$ws3.Range("N11").Value = -1676
# BSM row 12
$ws3.Range("H12").Value = 1788.5714
$ws3.Range("J12").Value = 3030
$ws3.Range("L12").Value = 3030
$ws3.Range("N12").Value = -3366
# BSM row 18
$ws3.Range("H18").Value = 8500
$ws3.Range("J18").Value = 8500
$ws3.Range("L18").Value = 8500
$ws3.Range("N18").Value = -9558
# BSM row 35
$ws3.Range("H35").Value = 60000
$ws3.Range("J35").Value = 60000
$ws3.Range("L35").Value = 60000
$ws3.Range("N35").Value = -60620
# BSM row 86
$ws3.Range("H86").Value = 1829.3235
$ws3.Range("I86").Value = 1272.8422
$ws3.Range("J86").Value = 2534.2
$ws3.Range("K86").Value = 1272.8422
$ws3.Range("L86").Value = 2534.2
$ws3.Range("M86").Value = -149.8422
$ws3.Range("N86").Value = -4780.2
# BSM row 89
$ws3.Range("H89").Value = 1829.3235
$ws3.Range("I89").Value = 1272.8422
$ws3.Range("J89").Value = 2534.2
$ws3.Range("K89").Value = 6364.211
$ws3.Range("L89").Value = 12671
$ws3.Range("M89").Value = -748.2110000000002
$ws3.Range("N89").Value = -23903
# BSM row 107
$ws3.Range("H107").Value = 166667840
$ws3.Range("I107").Value = 1000
$ws3.Range("J107").Value = 250001250
$ws3.Range("K107").Value = 1000
$ws3.Range("L107").Value = 250001250
$ws3.Range("M107").Value = 920
$ws3.Range("N107").Value = -250005090

$ws4 = $wb.Worksheets.Item("CRP")
# CRP row 11
$ws4.Range("H11").Value = 125001520
$ws4.Range("I11").Value = 300
$ws4.Range("J11").Value = 200002260
$ws4.Range("K11").Value = 300
$ws4.Range("L11").Value = 200002260
$ws4.Range("M11").Value = -160
$ws4.Range("N11").Value = -200002540
# CRP row 15
$ws4.Range("H15").Value = 0
$ws4.Range("I15").Value = 0
$ws4.Range("K15").Value = 0
$ws4.Range("M15").ClearContents()
# CRP row 31
$ws4.Range("H31").Value = 2766.689
$ws4.Range("I31").Value = 2002.65
$ws4.Range("J31").Value = 2984.9856
$ws4.Range("K31").Value = 2002.65
$ws4.Range("L31").Value = 2984.9856
$ws4.Range("M31").Value = -1707.65
$ws4.Range("N31").Value = -3574.9856
# CRP row 34
$ws4.Range("H34").Value = 2766.689
$ws4.Range("I34").Value = 2002.65
$ws4.Range("J34").Value = 2984.9856
$ws4.Range("K34").Value = 2002.65
$ws4.Range("L34").Value = 2984.9856
$ws4.Range("M34").Value = -1800.65
$ws4.Range("N34").Value = -3388.9856
# CRP row 62
$ws4.Range("H62").Value = 43536.08
$ws4.Range("I62").Value = 3808.625
$ws4.Range("J62").Value = 107100
$ws4.Range("K62").Value = 3808.625
$ws4.Range("L62").Value = 107100
$ws4.Range("M62").Value = -3184.625
$ws4.Range("N62").Value = -108348
# CRP row 65
$ws4.Range("H65").Value = 43536.08
$ws4.Range("I65").Value = 3808.625
$ws4.Range("J65").Value = 107100
$ws4.Range("K65").Value = 19043.125
$ws4.Range("L65").Value = 535500
$ws4.Range("M65").Value = -15923.125
$ws4.Range("N65").Value = -541740
# CRP row 105
$ws4.Range("H105").Value = 961.375
$ws4.Range("I105").Value = 948.3333
$ws4.Range("J105").Value = 1000.5
$ws4.Range("K105").Value = 948.3333
$ws4.Range("L105").Value = 1000.5
$ws4.Range("M105").Value = 798.6667
$ws4.Range("N105").Value = -4494.5
# CRP row 141
$ws4.Range("H141").Value = 290044.62
$ws4.Range("J141").Value = 290044.62
$ws4.Range("L141").Value = 290044.62
$ws4.Range("N141").Value = -300404.62

$ws6 = $wb.Worksheets.Item("GSM")
# GSM row 17
$ws6.Range("H17").Value = 425
$ws6.Range("I17").Value = 504
$ws6.Range("J17").Value = 372.33334
$ws6.Range("K17").Value = 504
$ws6.Range("L17").Value = 372.33334
$ws6.Range("M17").Value = -336
$ws6.Range("N17").Value = -708.33334
# GSM row 95
$ws6.Range("H95").Value = 11999.5
$ws6.Range("J95").Value = 11999.5
$ws6.Range("L95").Value = 11999.5
$ws6.Range("N95").Value = -17491.5
# GSM row 132
$ws6.Range("H132").Value = 4527.0586
$ws6.Range("I132").Value = 3707.8572
$ws6.Range("K132").Value = 11123.5716
$ws6.Range("M132").Value = -8593.571599999999

$ws7 = $wb.Worksheets.Item("LTW")
# LTW row 40
$ws7.Range("H40").Value = 2725.375
$ws7.Range("I40").Value = 2725.375
$ws7.Range("J40").Value = 0
$ws7.Range("K40").Value = 2725.375
$ws7.Range("L40").Value = 0
$ws7.Range("M40").Value = -2589.375
$ws7.Range("N40").ClearContents()
# LTW row 68
$ws7.Range("H68").Value = 28875
$ws7.Range("J68").Value = 28875
$ws7.Range("L68").Value = 28875
$ws7.Range("N68").Value = -30373
# LTW row 70
$ws7.Range("H70").Value = 19163
$ws7.Range("J70").Value = 19163
$ws7.Range("L70").Value = 19163
$ws7.Range("N70").Value = -19703
# LTW row 71
$ws7.Range("H71").Value = 28875
$ws7.Range("J71").Value = 28875
$ws7.Range("L71").Value = 144375
$ws7.Range("N71").Value = -151863
# LTW row 73
$ws7.Range("H73").Value = 19163
$ws7.Range("J73").Value = 19163
$ws7.Range("L73").Value = 19163
$ws7.Range("N73").Value = -21035
# LTW row 82
$ws7.Range("H82").Value = 1115.6666
$ws7.Range("I82").Value = 1148.7142
$ws7.Range("K82").Value = 1148.7142
$ws7.Range("M82").Value = -787.7141999999999
# LTW row 85
$ws7.Range("H85").Value = 1115.6666
$ws7.Range("I85").Value = 1148.7142
$ws7.Range("K85").Value = 1148.7142
$ws7.Range("M85").Value = 99.28580000000011

$ws8 = $wb.Worksheets.Item("WVR")
# WVR row 81
$ws8.Range("H81").Value = 10530317
$ws8.Range("I81").Value = 3113.6667
$ws8.Range("J81").Value = 20004800
$ws8.Range("K81").Value = 6227.3334
$ws8.Range("L81").Value = 40009600
$ws8.Range("M81").Value = -5166.3334
$ws8.Range("N81").Value = -40011722
# WVR row 84
$ws8.Range("H84").Value = 10530317
$ws8.Range("I84").Value = 3113.6667
$ws8.Range("J84").Value = 20004800
$ws8.Range("K84").Value = 31136.667
$ws8.Range("L84").Value = 200048000
$ws8.Range("M84").Value = -25832.667
$ws8.Range("N84").Value = -200058608
